# Update simulation output values (new input files re-run) for both worksheets.
# Only numeric result cells change; row/col structure and labels are untouched.
$wb = $excel.ActiveWorkbook
$wsOutput = $wb.Worksheets.Item("Output_flows")
$wsInput  = $wb.Worksheets.Item("Input_flows")

# --- Output_flows sheet ---
# Row 2
$wsOutput.Range("C2").Value = [double]"8.433726643683042E-33"
$wsOutput.Range("E2").Value = [double]"1.476432994525518E-28"
$wsOutput.Range("G2").Value = [double]"6.083647802527324E-30"
$wsOutput.Range("I2").Value = [double]"7.37620448702451E-30"
# Row 3
$wsOutput.Range("C3").Value = [double]"1.437235641732149E-33"
$wsOutput.Range("D3").Value = [double]"1.420200501997149E-33"
$wsOutput.Range("E3").Value = [double]"6.533906986096025E-28"
$wsOutput.Range("G3").Value = [double]"1.036746366457919E-30"
$wsOutput.Range("I3").Value = [double]"1.257017738107122E-30"
# Row 4
$wsOutput.Range("C4").Value = [double]"9.402449253602028E-34"
$wsOutput.Range("D4").Value = [double]"9.291004733138208E-33"
$wsOutput.Range("E4").Value = [double]"3.565757140670418E-25"
$wsOutput.Range("G4").Value = [double]"6.782433455190891E-31"
$wsOutput.Range("I4").Value = [double]"8.223456996366696E-31"
# Row 5
$wsOutput.Range("C5").Value = [double]"4.553004289116487E-34"
$wsOutput.Range("D5").Value = [double]"4.499038841818176E-32"
$wsOutput.Range("E5").Value = [double]"4.572930629786626E-22"
$wsOutput.Range("G5").Value = [double]"3.284298354527268E-31"
$wsOutput.Range("I5").Value = [double]"3.982093810448267E-31"
# Row 6
$wsOutput.Range("C6").Value = [double]"2.870745593686062E-32"
$wsOutput.Range("D6").Value = [double]"2.836719473742963E-29"
$wsOutput.Range("E6").Value = [double]"2.108163543598796E-16"
$wsOutput.Range("G6").Value = [double]"2.070805215832363E-29"
$wsOutput.Range("I6").Value = [double]"2.510776958263575E-29"
# Row 7
$wsOutput.Range("C7").Value = [double]"4.506114953801109E-25"
$wsOutput.Range("F7").Value = [double]"8.401842666699947E-29"
$wsOutput.Range("G7").Value = [double]"3.25047484876223E-21"
$wsOutput.Range("I7").Value = [double]"3.941084024364442E-21"
$wsOutput.Range("J7").Value = [double]"3.060702083515264E-21"
# Row 8
$wsOutput.Range("C8").Value = [double]"1.414553556389631E-25"
$wsOutput.Range("D8").Value = [double]"1.398253034736891E-26"
$wsOutput.Range("F8").Value = [double]"6.429418053567421E-28"
$wsOutput.Range("G8").Value = [double]"1.020384700437583E-21"
$wsOutput.Range("I8").Value = [double]"1.237179805631107E-21"
# Row 9
$wsOutput.Range("C9").Value = [double]"9.402409122495156E-26"
$wsOutput.Range("D9").Value = [double]"9.290968174681108E-26"
$wsOutput.Range("F9").Value = [double]"3.565707387818932E-25"
$wsOutput.Range("G9").Value = [double]"6.782404506716456E-22"
$wsOutput.Range("I9").Value = [double]"8.223421897381032E-22"
# Row 10
$wsOutput.Range("C10").Value = [double]"4.553004452355377E-26"
$wsOutput.Range("D10").Value = [double]"4.499039004621923E-25"
$wsOutput.Range("F10").Value = [double]"4.572930617265496E-22"
$wsOutput.Range("G10").Value = [double]"3.284298472279237E-22"
$wsOutput.Range("I10").Value = [double]"3.982093953218299E-22"
# Row 11
$wsOutput.Range("C11").Value = [double]"2.87074559383929E-24"
$wsOutput.Range("D11").Value = [double]"2.83671947389532E-22"
$wsOutput.Range("F11").Value = [double]"2.108163543599485E-16"
$wsOutput.Range("G11").Value = [double]"2.070805215942893E-20"
$wsOutput.Range("I11").Value = [double]"2.510776958397588E-20"
# Row 12
$wsOutput.Range("C12").Value = [double]"1.331363849080941E-33"
$wsOutput.Range("E12").Value = [double]"8.644964282618561E-26"
$wsOutput.Range("I12").Value = [double]"5.822107125566287E-30"
$wsOutput.Range("J12").Value = [double]"9.956069621431432E-28"
# Row 13
$wsOutput.Range("C13").Value = [double]"1.436845677993276E-34"
$wsOutput.Range("D13").Value = [double]"1.064861370289541E-33"
$wsOutput.Range("E13").Value = [double]"3.108887536930672E-26"
$wsOutput.Range("I13").Value = [double]"6.283383363577566E-31"
$wsOutput.Range("J13").Value = [double]"2.108104981612279E-28"
# Row 14
$wsOutput.Range("C14").Value = [double]"1.413718963806574E-35"
$wsOutput.Range("D14").Value = [double]"4.190887681426907E-34"
$wsOutput.Range("E14").Value = [double]"2.540189807777508E-25"
$wsOutput.Range("I14").Value = [double]"6.18224932155721E-32"
$wsOutput.Range("J14").Value = [double]"1.362290945511132E-28"
# Row 15
$wsOutput.Range("C15").Value = [double]"1.590368346253297E-34"
$wsOutput.Range("D15").Value = [double]"4.007371330227641E-32"
$wsOutput.Range("E15").Value = [double]"8.771088439911679E-22"
$wsOutput.Range("I15").Value = [double]"6.954744104992948E-31"
$wsOutput.Range("J15").Value = [double]"3.980664161168591E-27"
# Row 16
$wsOutput.Range("C16").Value = [double]"1.73529450218752E-30"
$wsOutput.Range("D16").Value = [double]"4.295390024543481E-27"
$wsOutput.Range("E16").Value = [double]"1.117503236891844E-13"
$wsOutput.Range("I16").Value = [double]"7.588511955703359E-27"
# Row 17
$wsOutput.Range("C17").Value = [double]"2.987329637901332E-25"
$wsOutput.Range("F17").Value = [double]"7.760667874548491E-26"
$wsOutput.Range("I17").Value = [double]"5.225484583571057E-21"
$wsOutput.Range("J17").Value = [double]"8.93797537360159E-19"
# Row 18
$wsOutput.Range("C18").Value = [double]"3.354521374263839E-26"
$wsOutput.Range("D18").Value = [double]"9.944406365787149E-27"
$wsOutput.Range("F18").Value = [double]"2.903487646019201E-26"
$wsOutput.Range("I18").Value = [double]"5.867782217294848E-22"
$wsOutput.Range("J18").Value = [double]"1.968839398484176E-19"
# Row 19
$wsOutput.Range("C19").Value = [double]"3.515839841579076E-27"
$wsOutput.Range("D19").Value = [double]"4.169001930785562E-27"
$wsOutput.Range("F19").Value = [double]"2.526931375027546E-25"
$wsOutput.Range("I19").Value = [double]"6.149963049736656E-23"
$wsOutput.Range("J19").Value = [double]"1.355181341749697E-19"
# Row 20
$wsOutput.Range("C20").Value = [double]"3.975929755170759E-26"
$wsOutput.Range("D20").Value = [double]"4.007380291342228E-25"
$wsOutput.Range("F20").Value = [double]"8.771108693984774E-22"
$wsOutput.Range("I20").Value = [double]"6.95475965471356E-22"
$wsOutput.Range("J20").Value = [double]"3.980673739426102E-18"
# Row 21
$wsOutput.Range("C21").Value = [double]"4.338236255732394E-22"
$wsOutput.Range("D21").Value = [double]"4.295390024805895E-20"
$wsOutput.Range("F21").Value = [double]"1.117503236891952E-13"
$wsOutput.Range("I21").Value = [double]"7.588511956164443E-18"

# --- Input_flows sheet ---
# Row 2
$wsInput.Range("C2").Value = [double]"6.871573293890415E-34"
# Row 22
$wsInput.Range("E22").Value = [double]"2.096910438516114E-29"
# Row 23
$wsInput.Range("E23").Value = [double]"3.172366062988177E-30"
# Row 24
$wsInput.Range("E24").Value = [double]"1.380861982386761E-30"
# Row 25
$wsInput.Range("E25").Value = [double]"2.291897590323544E-31"
# Row 26
$wsInput.Range("E26").Value = [double]"6.313912933249146E-31"
# Row 27
$wsInput.Range("E27").Value = [double]"3.811874689118481E-21"
# Row 28
$wsInput.Range("E28").Value = [double]"8.385549295963613E-22"
# Row 29
$wsInput.Range("E29").Value = [double]"5.573746730425661E-22"
# Row 30
$wsInput.Range("E30").Value = [double]"2.698745507937039E-22"
# Row 31
$wsInput.Range("E31").Value = [double]"1.615518470865224E-20"
# Row 32
$wsInput.Range("E32").Value = [double]"2.59463656317397E-27"
# Row 33
$wsInput.Range("E33").Value = [double]"6.153469801912532E-28"
# Row 34
$wsInput.Range("E34").Value = [double]"3.919523913078316E-28"
# Row 35
$wsInput.Range("E35").Value = [double]"5.72638436643636E-28"
# Row 36
$wsInput.Range("E36").Value = [double]"1.776100409646976E-28"
# Row 37
$wsInput.Range("E37").Value = [double]"3.288176281641159E-19"
# Row 38
$wsInput.Range("E38").Value = [double]"7.217031566299576E-20"
# Row 39
$wsInput.Range("E39").Value = [double]"4.96082939467926E-20"
# Row 40
$wsInput.Range("E40").Value = [double]"1.422700202928972E-18"
# Row 41
$wsInput.Range("E41").Value = [double]"2.695371803750808E-18"
# Row 42
$wsInput.Range("E42").Value = [double]"3.313267739150016E-29"
# Row 43
$wsInput.Range("E43").Value = [double]"5.012564170700179E-30"
# Row 44
$wsInput.Range("E44").Value = [double]"2.181860214162712E-30"
# Row 45
$wsInput.Range("E45").Value = [double]"3.621361317094852E-31"
# Row 46
$wsInput.Range("E46").Value = [double]"9.976431823354519E-31"
# Row 47
$wsInput.Range("E47").Value = [double]"6.023033316614253E-21"
# Row 48
$wsInput.Range("E48").Value = [double]"1.324976472387133E-21"
# Row 49
$wsInput.Range("E49").Value = [double]"8.806916542024883E-22"
# Row 50
$wsInput.Range("E50").Value = [double]"4.264209984071335E-22"
# Row 51
$wsInput.Range("E51").Value = [double]"2.552634167488112E-20"
# Row 52
$wsInput.Range("E52").Value = [double]"4.09971044145638E-27"
# Row 53
$wsInput.Range("E53").Value = [double]"9.722920256402719E-28"
# Row 54
$wsInput.Range("E54").Value = [double]"6.193126752337267E-28"
# Row 55
$wsInput.Range("E55").Value = [double]"9.048094870810476E-28"
# Row 56
$wsInput.Range("E56").Value = [double]"2.806365060082766E-28"
# Row 57
$wsInput.Range("E57").Value = [double]"5.19555255889208E-19"
# Row 58
$wsInput.Range("E58").Value = [double]"1.140342354248047E-19"
# Row 59
$wsInput.Range("E59").Value = [double]"7.838463527535912E-20"
# Row 60
$wsInput.Range("E60").Value = [double]"2.247967580428696E-18"
# Row 61
$wsInput.Range("E61").Value = [double]"4.258879291335796E-18"
# Row 62
$wsInput.Range("E62").Value = [double]"1.104422579716672E-30"
# Row 63
$wsInput.Range("E63").Value = [double]"1.670854723566727E-31"
# Row 64
$wsInput.Range("E64").Value = [double]"7.272867380542371E-32"
# Row 65
$wsInput.Range("E65").Value = [double]"1.207120439031618E-32"
# Row 66
$wsInput.Range("E66").Value = [double]"3.325477274451507E-32"
# Row 67
$wsInput.Range("E67").Value = [double]"2.007677772204751E-22"
# Row 68
$wsInput.Range("E68").Value = [double]"4.416588241290444E-23"
# Row 69
$wsInput.Range("E69").Value = [double]"2.935638847341629E-23"
# Row 70
$wsInput.Range("E70").Value = [double]"1.421403328023779E-23"
# Row 71
$wsInput.Range("E71").Value = [double]"8.508780558293708E-22"
# Row 72
$wsInput.Range("E72").Value = [double]"1.366570147152127E-28"
# Row 73
$wsInput.Range("E73").Value = [double]"3.240973418800907E-29"
# Row 74
$wsInput.Range("E74").Value = [double]"2.064375584112423E-29"
# Row 75
$wsInput.Range("E75").Value = [double]"3.016031623603492E-29"
# Row 76
$wsInput.Range("E76").Value = [double]"9.354550200275888E-30"
# Row 77
$wsInput.Range("E77").Value = [double]"1.731850852964027E-20"
# Row 78
$wsInput.Range("E78").Value = [double]"3.801141180826823E-21"
# Row 79
$wsInput.Range("E79").Value = [double]"2.612821175845304E-21"
# Row 80
$wsInput.Range("E80").Value = [double]"7.493225268095654E-20"
# Row 81
$wsInput.Range("E81").Value = [double]"1.419626430445265E-19"
# Row 82
$wsInput.Range("F82").Value = [double]"1.996171035920316E-30"
# Row 83
$wsInput.Range("F83").Value = [double]"1.247762968593766E-31"
# Row 84
$wsInput.Range("F84").Value = [double]"4.937021608066122E-33"
# Row 85
$wsInput.Range("F85").Value = [double]"2.441052400954396E-33"
# Row 86
$wsInput.Range("F86").Value = [double]"3.573631359695709E-30"
# Row 87
$wsInput.Range("F87").Value = [double]"1.813497867506498E-23"
# Row 88
$wsInput.Range("F88").Value = [double]"7.335421379988738E-24"
# Row 89
$wsInput.Range("F89").Value = [double]"4.889677969918797E-24"
# Row 90
$wsInput.Range("F90").Value = [double]"2.441495392074335E-24"
# Row 91
$wsInput.Range("F91").Value = [double]"3.56996006601156E-21"
# Row 92
$wsInput.Range("F92").Value = [double]"5.084783709905169E-31"
# Row 93
$wsInput.Range("F93").Value = [double]"3.673980366437109E-32"
# Row 94
$wsInput.Range("F94").Value = [double]"5.423411544436765E-34"
# Row 95
$wsInput.Range("F95").Value = [double]"2.332863876995667E-28"
# Row 96
$wsInput.Range("F96").Value = [double]"5.149622272033073E-28"
# Row 97
$wsInput.Range("F97").Value = [double]"1.346175663746884E-23"
# Row 98
$wsInput.Range("F98").Value = [double]"1.515267387754151E-24"
# Row 99
$wsInput.Range("F99").Value = [double]"1.791433851588627E-25"
# Row 100
$wsInput.Range("F100").Value = [double]"2.332934970160521E-19"
# Row 101
$wsInput.Range("F101").Value = [double]"5.149779006331733E-19"
